$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5028.1
$ws.Range("I19").Value = 1744.7307
$ws.Range("K19").Value = 1744.7307
$ws.Range("M19").Value = -1569.7307
$ws.Range("H40").Value = 7116
$ws.Range("I40").Value = 6563.6665
$ws.Range("J40").Value = 8496.833000000001
$ws.Range("K40").Value = 6563.6665
$ws.Range("L40").Value = 8496.833000000001
$ws.Range("M40").Value = -6388.6665
$ws.Range("N40").Value = -8846.833000000001
$ws.Range("H94").Value = 8419323
$ws.Range("I94").Value = 10223003
$ws.Range("J94").Value = 2151.6667
$ws.Range("K94").Value = 10223003
$ws.Range("L94").Value = 2151.6667
$ws.Range("M94").Value = -10222552
$ws.Range("N94").Value = -3053.6667
$ws.Range("H111").Value = 932.2778
$ws.Range("I111").Value = 878.25
$ws.Range("J111").Value = 1040.3334
$ws.Range("K111").Value = 2634.75
$ws.Range("L111").Value = 3121.0002
$ws.Range("M111").Value = 432.25
$ws.Range("N111").Value = -9255.0002
$ws.Range("H112").Value = 3036.6562
$ws.Range("J112").Value = 3153.0715
$ws.Range("L112").Value = 9459.2145
$ws.Range("N112").Value = -11675.2145
$ws.Range("H125").Value = 4770.048
$ws.Range("I125").Value = 1675
$ws.Range("J125").Value = 9799.5
$ws.Range("K125").Value = 15075
$ws.Range("L125").Value = 88195.5
$ws.Range("M125").Value = -12615
$ws.Range("N125").Value = -93115.5
$ws.Range("H137").Value = 3720.0188
$ws.Range("I137").Value = 4441.946
$ws.Range("J137").Value = 2050.5625
$ws.Range("K137").Value = 13325.838
$ws.Range("L137").Value = 6151.6875
$ws.Range("M137").Value = -10775.838
$ws.Range("N137").Value = -11251.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4480
$ws.Range("H40").Value = 36999.668
$ws.Range("I40").Value = 34499.5
$ws.Range("J40").Value = 42000
$ws.Range("K40").Value = 34499.5
$ws.Range("L40").Value = 42000
$ws.Range("M40").Value = -34323.5
$ws.Range("N40").Value = -42352
$ws.Range("H61").Value = 4261.104
$ws.Range("I61").Value = 3813.5083
$ws.Range("J61").Value = 5967.5625
$ws.Range("K61").Value = 3813.5083
$ws.Range("L61").Value = 5967.5625
$ws.Range("M61").Value = -3601.5083
$ws.Range("N61").Value = -6391.5625
$ws.Range("H74").Value = 2366.4849
$ws.Range("I74").Value = 1842.5454
$ws.Range("J74").Value = 3414.3635
$ws.Range("K74").Value = 1842.5454
$ws.Range("L74").Value = 3414.3635
$ws.Range("M74").Value = -968.5454
$ws.Range("N74").Value = -5162.363499999999
$ws.Range("H77").Value = 2366.4849
$ws.Range("I77").Value = 1842.5454
$ws.Range("J77").Value = 3414.3635
$ws.Range("K77").Value = 9212.726999999999
$ws.Range("L77").Value = 17071.8175
$ws.Range("M77").Value = -4844.726999999999
$ws.Range("N77").Value = -25807.8175
$ws.Range("H102").Value = 3132.4546
$ws.Range("I102").Value = 1065.8572
$ws.Range("K102").Value = 1065.8572
$ws.Range("M102").Value = 556.1428000000001
$ws.Range("H136").Value = 4261.104
$ws.Range("I136").Value = 3813.5083
$ws.Range("J136").Value = 5967.5625
$ws.Range("K136").Value = 11440.5249
$ws.Range("L136").Value = 17902.6875
$ws.Range("M136").Value = -8890.5249
$ws.Range("N136").Value = -23002.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3690.6191
$ws.Range("I20").Value = 3579.2856
$ws.Range("J20").Value = 3913.2856
$ws.Range("K20").Value = 3579.2856
$ws.Range("L20").Value = 3913.2856
$ws.Range("M20").Value = -3332.2856
$ws.Range("N20").Value = -4407.2856
$ws.Range("H22").Value = 579.25
$ws.Range("I22").Value = 594.7368
$ws.Range("K22").Value = 594.7368
$ws.Range("M22").Value = -421.7368
$ws.Range("H99").Value = 4899.3335
$ws.Range("I99").Value = 4279.2
$ws.Range("K99").Value = 4279.2
$ws.Range("M99").Value = -2781.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 5418.625
$ws.Range("I33").Value = 4659.8
$ws.Range("J33").Value = 6683.3335
$ws.Range("K33").Value = 4659.8
$ws.Range("L33").Value = 6683.3335
$ws.Range("M33").Value = -4280.8
$ws.Range("N33").Value = -7441.3335
$ws.Range("H42").Value = 26264
$ws.Range("I42").Value = 28352
$ws.Range("J42").Value = 20000
$ws.Range("K42").Value = 28352
$ws.Range("L42").Value = 20000
$ws.Range("M42").Value = -27759
$ws.Range("N42").Value = -21186
$ws.Range("H58").Value = 4946
$ws.Range("I58").Value = 7000
$ws.Range("J58").Value = 3919
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 3919
$ws.Range("M58").Value = -6797
$ws.Range("N58").Value = -4325
$ws.Range("H125").Value = 90000
$ws.Range("J125").Value = 90000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920
$ws.Range("H136").Value = 4946
$ws.Range("I136").Value = 7000
$ws.Range("J136").Value = 3919
$ws.Range("K136").Value = 21000
$ws.Range("L136").Value = 11757
$ws.Range("M136").Value = -18450
$ws.Range("N136").Value = -16857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 28.461538
$ws.Range("I12").Value = 2.3333333
$ws.Range("K12").Value = 6.999999900000001
$ws.Range("M12").Value = 166.0000001
$ws.Range("H61").Value = 193.3077
$ws.Range("I61").Value = 188.54546
$ws.Range("K61").Value = 565.6363799999999
$ws.Range("M61").Value = -350.6363799999999
$ws.Range("H75").Value = 1015
$ws.Range("J75").Value = 1015
$ws.Range("L75").Value = 3045
$ws.Range("N75").Value = -5041
$ws.Range("H78").Value = 1015
$ws.Range("J78").Value = 1015
$ws.Range("L78").Value = 9135
$ws.Range("N78").Value = -19119
$ws.Range("H103").Value = 814.875
$ws.Range("I103").Value = 878.3333
$ws.Range("J103").Value = 624.5
$ws.Range("K103").Value = 2634.9999
$ws.Range("L103").Value = 1873.5
$ws.Range("M103").Value = -1755.9999
$ws.Range("N103").Value = -3631.5
$ws.Range("H124").Value = 7997
$ws.Range("I124").Value = 3494
$ws.Range("K124").Value = 10482
$ws.Range("M124").Value = -5572
$ws.Range("H130").Value = 5406.3335
$ws.Range("I130").Value = 1859.5
$ws.Range("K130").Value = 5578.5
$ws.Range("M130").Value = -558.5
$ws.Range("H131").Value = 701335.8
$ws.Range("I131").Value = 919893.25
$ws.Range("J131").Value = 1952.1
$ws.Range("K131").Value = 2759679.75
$ws.Range("L131").Value = 5856.299999999999
$ws.Range("M131").Value = -2754639.75
$ws.Range("N131").Value = -15936.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 26766.334
$ws.Range("J18").Value = 40400
$ws.Range("L18").Value = 40400
$ws.Range("N18").Value = -40986
$ws.Range("H35").Value = 34498.5
$ws.Range("I35").Value = 25000
$ws.Range("J35").Value = 43997
$ws.Range("K35").Value = 25000
$ws.Range("L35").Value = 43997
$ws.Range("M35").Value = -24702
$ws.Range("N35").Value = -44593
$ws.Range("H42").Value = 130287.664
$ws.Range("J42").Value = 130287.664
$ws.Range("L42").Value = 130287.664
$ws.Range("N42").Value = -131257.664
$ws.Range("H102").Value = 6742.982
$ws.Range("I102").Value = 5235
$ws.Range("J102").Value = 11159.214
$ws.Range("K102").Value = 5235
$ws.Range("L102").Value = 11159.214
$ws.Range("M102").Value = -3613
$ws.Range("N102").Value = -14403.214
$ws.Range("H115").Value = 130287.664
$ws.Range("J115").Value = 130287.664
$ws.Range("L115").Value = 130287.664
$ws.Range("N115").Value = -132637.664
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 19286.715
$ws.Range("J20").Value = 20333.666
$ws.Range("L20").Value = 20333.666
$ws.Range("N20").Value = -20785.666
$ws.Range("H22").Value = 824
$ws.Range("I22").Value = 868.5833
$ws.Range("J22").Value = 734.8333
$ws.Range("K22").Value = 868.5833
$ws.Range("L22").Value = 734.8333
$ws.Range("M22").Value = -573.5833
$ws.Range("N22").Value = -1324.8333
$ws.Range("H27").Value = 824
$ws.Range("I27").Value = 868.5833
$ws.Range("J27").Value = 734.8333
$ws.Range("K27").Value = 868.5833
$ws.Range("L27").Value = 734.8333
$ws.Range("M27").Value = -761.5833
$ws.Range("N27").Value = -948.8333
$ws.Range("H46").Value = 2881.7222
$ws.Range("I46").Value = 963.3333
$ws.Range("K46").Value = 963.3333
$ws.Range("M46").Value = -775.3333
$ws.Range("H55").Value = 767.875
$ws.Range("I55").Value = 726.3333
$ws.Range("J55").Value = 821.2857
$ws.Range("K55").Value = 726.3333
$ws.Range("L55").Value = 821.2857
$ws.Range("M55").Value = -553.3333
$ws.Range("N55").Value = -1167.2857
$ws.Range("H132").Value = 23937.797
$ws.Range("I132").Value = 36046.656
$ws.Range("K132").Value = 108139.968
$ws.Range("M132").Value = -105609.968

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10478.3
$ws.Range("I41").Value = 2250
$ws.Range("J41").Value = 12535.375
$ws.Range("K41").Value = 2250
$ws.Range("L41").Value = 12535.375
$ws.Range("M41").Value = -1860
$ws.Range("N41").Value = -13315.375
$ws.Range("H118").Value = 104999
$ws.Range("J118").Value = 104999
$ws.Range("L118").Value = 104999
$ws.Range("N118").Value = -108313
$ws.Range("H128").Value = 59724
$ws.Range("J128").Value = 59724
$ws.Range("L128").Value = 59724
$ws.Range("N128").Value = -69684
$ws.Range("H132").Value = 5651.033
$ws.Range("I132").Value = 4008.976
$ws.Range("J132").Value = 9482.5
$ws.Range("K132").Value = 12026.928
$ws.Range("L132").Value = 28447.5
$ws.Range("M132").Value = -9496.928
$ws.Range("N132").Value = -33507.5
